$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 324 (old row 324 becomes row 325), shifting cells down.
$ws.Rows.Item(324).Insert(-4121)

# Populate the newly inserted row 324 with a duplicate of the (now shifted) row 323 data,
# matching the target state described in the diff.
$ws.Range("A324").Value = 8
$ws.Range("B324").Value = "Terminal La Palmera de La Serena"
$ws.Range("C324").Value = "Coquimbo"
$ws.Range("D324").Value = 44736
$ws.Range("E324").Value = 4
$ws.Range("F324").Value = 100112031
$ws.Range("G324").Value = "Poroto verde"
$ws.Range("H324").Value = "Magnum"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 400
$ws.Range("K324").Value = 25000
$ws.Range("L324").Value = 26000
$ws.Range("M324").Value = 25500
$ws.Range("N324").Value = "$/malla 25 kilos"
$ws.Range("O324").Value = "Perú"
$ws.Range("P324").Value = 1020
$ws.Range("Q324").Value = 25
$ws.Range("R324").Value = "Hortaliza"
